$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 189
$ws.Cells.Item(189, 8).Value = 0
$ws.Cells.Item(189, 9).Value = 3
$ws.Cells.Item(189, 10).Value = "A"
$ws.Cells.Item(189, 14).Value = 2.15
$ws.Cells.Item(189, 15).Value = 3.2
$ws.Cells.Item(189, 16).Value = 3.2
$ws.Cells.Item(189, 17).Value = -0.25
$ws.Cells.Item(189, 18).Value = 1.9
$ws.Cells.Item(189, 19).Value = 1.9
$ws.Cells.Item(189, 21).Value = 1.95
$ws.Cells.Item(189, 22).Value = 1.85
$ws.Cells.Item(189, 23).Value = -1
$ws.Cells.Item(189, 24).Value = -1
$ws.Cells.Item(189, 25).Value = 2.2
$ws.Cells.Item(189, 26).Value = -1
$ws.Cells.Item(189, 27).Value = 0.8999999999999999
$ws.Cells.Item(189, 28).Value = 0.95
$ws.Cells.Item(189, 29).Value = -1

# Row 190
$ws.Cells.Item(190, 8).Value = 2
$ws.Cells.Item(190, 9).Value = 1
$ws.Cells.Item(190, 10).Value = "H"
$ws.Cells.Item(190, 14).Value = 1.571
$ws.Cells.Item(190, 15).Value = 4.2
$ws.Cells.Item(190, 16).Value = 5.75
$ws.Cells.Item(190, 18).Value = 1.925
$ws.Cells.Item(190, 19).Value = 1.875
$ws.Cells.Item(190, 21).Value = 1.85
$ws.Cells.Item(190, 22).Value = 1.95
$ws.Cells.Item(190, 23).Value = 0.571
$ws.Cells.Item(190, 24).Value = -1
$ws.Cells.Item(190, 25).Value = -1
$ws.Cells.Item(190, 27).Value = -0
$ws.Cells.Item(190, 28).Value = 0.425
$ws.Cells.Item(190, 29).Value = -0.5

# Row 191
$ws.Cells.Item(191, 8).Value = 1
$ws.Cells.Item(191, 9).Value = 2
$ws.Cells.Item(191, 10).Value = "A"
$ws.Cells.Item(191, 14).Value = 6.5
$ws.Cells.Item(191, 15).Value = 3.6
$ws.Cells.Item(191, 16).Value = 1.6
$ws.Cells.Item(191, 17).Value = 1
$ws.Cells.Item(191, 18).Value = 1.725
$ws.Cells.Item(191, 19).Value = 1.975
$ws.Cells.Item(191, 20).Value = 2
$ws.Cells.Item(191, 21).Value = 1.775
$ws.Cells.Item(191, 22).Value = 2.025
$ws.Cells.Item(191, 23).Value = -1
$ws.Cells.Item(191, 24).Value = -1
$ws.Cells.Item(191, 25).Value = 0.6000000000000001
$ws.Cells.Item(191, 27).Value = -0
$ws.Cells.Item(191, 28).Value = 0.7749999999999999
$ws.Cells.Item(191, 29).Value = -1

# Row 192
$ws.Cells.Item(192, 8).Value = 2
$ws.Cells.Item(192, 9).Value = 2
$ws.Cells.Item(192, 10).Value = "D"
$ws.Cells.Item(192, 14).Value = 3
$ws.Cells.Item(192, 15).Value = 3.6
$ws.Cells.Item(192, 16).Value = 2.3
$ws.Cells.Item(192, 17).Value = 0.25
$ws.Cells.Item(192, 23).Value = -1
$ws.Cells.Item(192, 24).Value = 2.6
$ws.Cells.Item(192, 25).Value = -1
$ws.Cells.Item(192, 26).Value = 0.3875
$ws.Cells.Item(192, 27).Value = -0.5
$ws.Cells.Item(192, 28).Value = 0.95
$ws.Cells.Item(192, 29).Value = -1

# Row 193
$ws.Cells.Item(193, 8).Value = 0
$ws.Cells.Item(193, 9).Value = 1
$ws.Cells.Item(193, 10).Value = "A"
$ws.Cells.Item(193, 14).Value = 2.4
$ws.Cells.Item(193, 15).Value = 3.6
$ws.Cells.Item(193, 16).Value = 2.8
$ws.Cells.Item(193, 17).Value = 0
$ws.Cells.Item(193, 18).Value = 1.775
$ws.Cells.Item(193, 19).Value = 2.025
$ws.Cells.Item(193, 20).Value = 2.25
$ws.Cells.Item(193, 21).Value = 1.775
$ws.Cells.Item(193, 22).Value = 2.025
$ws.Cells.Item(193, 23).Value = -1
$ws.Cells.Item(193, 24).Value = -1
$ws.Cells.Item(193, 25).Value = 1.8
$ws.Cells.Item(193, 26).Value = -1
$ws.Cells.Item(193, 27).Value = 1.025
$ws.Cells.Item(193, 28).Value = -1
$ws.Cells.Item(193, 29).Value = 1.025

# Row 194
$ws.Cells.Item(194, 8).Value = 0
$ws.Cells.Item(194, 9).Value = 0
$ws.Cells.Item(194, 10).Value = "D"
$ws.Cells.Item(194, 14).Value = 3
$ws.Cells.Item(194, 15).Value = 3.2
$ws.Cells.Item(194, 16).Value = 2.5
$ws.Cells.Item(194, 17).Value = 0.25
$ws.Cells.Item(194, 18).Value = 1.75
$ws.Cells.Item(194, 19).Value = 2.05
$ws.Cells.Item(194, 21).Value = 2.025
$ws.Cells.Item(194, 22).Value = 1.775
$ws.Cells.Item(194, 23).Value = -1
$ws.Cells.Item(194, 24).Value = 2.2
$ws.Cells.Item(194, 25).Value = -1
$ws.Cells.Item(194, 26).Value = 0.375
$ws.Cells.Item(194, 27).Value = -0.5
$ws.Cells.Item(194, 28).Value = -1
$ws.Cells.Item(194, 29).Value = 0.7749999999999999

# Row 195
$ws.Cells.Item(195, 15).Value = 3.4
$ws.Cells.Item(195, 16).Value = 3.1
$ws.Cells.Item(195, 18).Value = 2.05
$ws.Cells.Item(195, 19).Value = 1.75
